$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (subject "002") and shift the existing subjects down.
$ws.Rows.Item(3).Insert()

# Copy the formatting from the row directly below (the former row 3, now row 4,
# which already carries the correct style set for a "standalone number" subject)
# onto the newly inserted, blank row 3.
$ws.Range("B4:G4").Copy()
$ws.Range("B3:G3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new subject's data.
$ws.Range("B3").Value = "002"
$ws.Range("C3").Value = 78
$ws.Range("D3").Value = "F"
$ws.Range("E3").Value = 75
$ws.Range("F3").Value = 1.62

# Extend the autofilter / sort range to cover the new last row (28).
$ws.Range("B2:G28").AutoFilter(1)
